$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 26584
$ws.Range("E2").Value = 518165158028
$ws.Range("F2").Value = 5839501584
$ws.Range("G2").Value = 0.28158

$ws.Range("D3").Value = 1631.21
$ws.Range("E3").Value = 196052662000
$ws.Range("F3").Value = 3906952610
$ws.Range("G3").Value = -0.21511

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 83065234589
$ws.Range("F4").Value = 10276939792
$ws.Range("G4").Value = 0.00714

$ws.Range("D5").Value = 214.91
$ws.Range("E5").Value = 33041536843
$ws.Range("F5").Value = 316341184
$ws.Range("G5").Value = 0.3917

$ws.Range("D6").Value = 0.49768
$ws.Range("E6").Value = 26465863331
$ws.Range("F6").Value = 372474800
$ws.Range("G6").Value = -0.6498

$ws.Range("D7").Value = 0.9999130000000001
$ws.Range("E7").Value = 26160397739
$ws.Range("F7").Value = 2066493106
$ws.Range("G7").Value = 0.00254

$ws.Range("D8").Value = 1631.47
$ws.Range("E8").Value = 14106344663
$ws.Range("F8").Value = 5269809
$ws.Range("G8").Value = -0.2427

$ws.Range("B9").Value = "DOGE"
$ws.Range("C9").Value = "Dogecoin"
$ws.Range("D9").Value = 0.062321
$ws.Range("E9").Value = 8789177308
$ws.Range("F9").Value = 167179588
$ws.Range("G9").Value = -0.05315

$ws.Range("B10").Value = "ADA"
$ws.Range("C10").Value = "Cardano"
$ws.Range("D10").Value = 0.250248
$ws.Range("E10").Value = 8767878478
$ws.Range("F10").Value = 73026140
$ws.Range("G10").Value = -0.30638

$ws.Range("B11").Value = "TON"
$ws.Range("C11").Value = "Toncoin"
$ws.Range("D11").Value = 2.43
$ws.Range("E11").Value = 8375803632
$ws.Range("F11").Value = 68178591
$ws.Range("G11").Value = 7.9063

$ws.Range("B12").Value = "SOL"
$ws.Range("C12").Value = "Solana"
$ws.Range("D12").Value = 19.06
$ws.Range("E12").Value = 7838277357
$ws.Range("F12").Value = 122627495
$ws.Range("G12").Value = -0.50687

$ws.Range("B13").Value = "TRX"
$ws.Range("C13").Value = "TRON"
$ws.Range("D13").Value = 0.083727
$ws.Range("E13").Value = 7462129731
$ws.Range("F13").Value = 176537629
$ws.Range("G13").Value = 0.1771

$ws.Range("D14").Value = 4.12
$ws.Range("E14").Value = 5257682319
$ws.Range("F14").Value = 67246162
$ws.Range("G14").Value = -1.02201

$ws.Range("D15").Value = 0.524705
$ws.Range("E15").Value = 4889796203
$ws.Range("F15").Value = 111531762
$ws.Range("G15").Value = -0.45731

$ws.Range("D16").Value = 64.19
$ws.Range("E16").Value = 4727271061
$ws.Range("F16").Value = 279974983
$ws.Range("G16").Value = -1.14495

$ws.Range("B17").Value = "WBTC"
$ws.Range("C17").Value = "Wrapped Bitcoin"
$ws.Range("D17").Value = 26607
$ws.Range("E17").Value = 4331716641
$ws.Range("F17").Value = 22514198
$ws.Range("G17").Value = 0.29601

$ws.Range("B18").Value = "SHIB"
$ws.Range("C18").Value = "Shiba Inu"
$ws.Range("D18").Value = 0.0000073
$ws.Range("E18").Value = 4304061575
$ws.Range("F18").Value = 59269706
$ws.Range("G18").Value = -1.634

$ws.Range("B19").Value = "BCH"
$ws.Range("C19").Value = "Bitcoin Cash"
$ws.Range("D19").Value = 212.73
$ws.Range("E19").Value = 4151114857
$ws.Range("F19").Value = 125780347
$ws.Range("G19").Value = -1.19077

$ws.Range("B20").Value = "DAI"
$ws.Range("C20").Value = "Dai"
$ws.Range("D20").Value = 0.999627
$ws.Range("E20").Value = 3842286310
$ws.Range("F20").Value = 40100406
$ws.Range("G20").Value = 0.09317

$ws.Range("B21").Value = "LEO"
$ws.Range("C21").Value = "LEO Token"
$ws.Range("D21").Value = 3.69
$ws.Range("E21").Value = 3437538747
$ws.Range("F21").Value = 227597
$ws.Range("G21").Value = 1.2463

$ws.Range("B22").Value = "LINK"
$ws.Range("C22").Value = "Chainlink"
$ws.Range("D22").Value = 6.21
$ws.Range("E22").Value = 3342759802
$ws.Range("F22").Value = 102251195
$ws.Range("G22").Value = -0.40433

$ws.Range("B23").Value = "UNI"
$ws.Range("C23").Value = "Uniswap"
$ws.Range("D23").Value = 4.36
$ws.Range("E23").Value = 3283186735
$ws.Range("F23").Value = 66321209
$ws.Range("G23").Value = 0.79449

$ws.Range("B24").Value = "AVAX"
$ws.Range("C24").Value = "Avalanche"
$ws.Range("D24").Value = 9.27
$ws.Range("E24").Value = 3280165452
$ws.Range("F24").Value = 78287052
$ws.Range("G24").Value = -2.0715

$ws.Range("B25").Value = "XLM"
$ws.Range("C25").Value = "Stellar"
$ws.Range("D25").Value = 0.117484
$ws.Range("E25").Value = 3237432412
$ws.Range("F25").Value = 43956628
$ws.Range("G25").Value = -1.34729

$ws.Range("B26").Value = "TUSD"
$ws.Range("C26").Value = "TrueUSD"
$ws.Range("D26").Value = 0.999344
$ws.Range("E26").Value = 3101719109
$ws.Range("F26").Value = 136856015
$ws.Range("G26").Value = 0.08

$ws.Range("B27").Value = "XMR"
$ws.Range("C27").Value = "Monero"
$ws.Range("D27").Value = 145.22
$ws.Range("E27").Value = 2636840774
$ws.Range("F27").Value = 32891264
$ws.Range("G27").Value = -0.16323

$ws.Range("B28").Value = "OKB"
$ws.Range("C28").Value = "OKB"
$ws.Range("D28").Value = 43.45
$ws.Range("E28").Value = 2606526110
$ws.Range("F28").Value = 6564723
$ws.Range("G28").Value = 1.55008

$ws.Range("B29").Value = "BUSD"
$ws.Range("C29").Value = "Binance USD"
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 2492512278
$ws.Range("F29").Value = 1348891962
$ws.Range("G29").Value = 0.0431

$ws.Range("D30").Value = 15.52
$ws.Range("E30").Value = 2219503426
$ws.Range("F30").Value = 37420130
$ws.Range("G30").Value = -0.62252

$ws.Range("D31").Value = 7.13
$ws.Range("E31").Value = 2084088817
$ws.Range("F31").Value = 87593162
$ws.Range("G31").Value = 0.01459

$ws.Range("D32").Value = 0.050395
$ws.Range("E32").Value = 1673977779
$ws.Range("F32").Value = 17388285
$ws.Range("G32").Value = -2.03297

$ws.Range("B33").Value = "FIL"
$ws.Range("C33").Value = "Filecoin"
$ws.Range("D33").Value = 3.33
$ws.Range("E33").Value = 1495739765
$ws.Range("F33").Value = 64115371
$ws.Range("G33").Value = -1.16372

$ws.Range("B34").Value = "CRO"
$ws.Range("C34").Value = "Cronos"
$ws.Range("D34").Value = 0.051696
$ws.Range("E34").Value = 1358552943
$ws.Range("F34").Value = 9605855
$ws.Range("G34").Value = 0.9466599999999999

$ws.Range("B35").Value = "LDO"
$ws.Range("C35").Value = "Lido DAO"
$ws.Range("D35").Value = 1.52
$ws.Range("E35").Value = 1353083285
$ws.Range("F35").Value = 21451619
$ws.Range("G35").Value = -0.56306

$ws.Range("B36").Value = "ICP"
$ws.Range("C36").Value = "Internet Computer"
$ws.Range("D36").Value = 2.99
$ws.Range("E36").Value = 1329876596
$ws.Range("F36").Value = 14847374
$ws.Range("G36").Value = -1.45445

$ws.Range("B37").Value = "QNT"
$ws.Range("C37").Value = "Quant"
$ws.Range("D37").Value = 91.36
$ws.Range("E37").Value = 1328352347
$ws.Range("F37").Value = 16181164
$ws.Range("G37").Value = 0.16614

$ws.Range("B38").Value = "MNT"
$ws.Range("C38").Value = "Mantle"
$ws.Range("D38").Value = 0.404231
$ws.Range("E38").Value = 1307214768
$ws.Range("F38").Value = 14721043
$ws.Range("G38").Value = -0.09135

$ws.Range("B39").Value = "VET"
$ws.Range("C39").Value = "VeChain"
$ws.Range("D39").Value = 0.01765241
$ws.Range("E39").Value = 1279049631
$ws.Range("F39").Value = 39276705
$ws.Range("G39").Value = -1.29015

$ws.Range("B40").Value = "APT"
$ws.Range("C40").Value = "Aptos"
$ws.Range("D40").Value = 5.3
$ws.Range("E40").Value = 1247192374
$ws.Range("F40").Value = 41858751
$ws.Range("G40").Value = -0.6208900000000001

$ws.Range("B41").Value = "MKR"
$ws.Range("C41").Value = "Maker"
$ws.Range("D41").Value = 1275.77
$ws.Range("E41").Value = 1149935656
$ws.Range("F41").Value = 62288340
$ws.Range("G41").Value = 0.77069

$ws.Range("B42").Value = "OP"
$ws.Range("C42").Value = "Optimism"
$ws.Range("D42").Value = 1.38
$ws.Range("E42").Value = 1097406751
$ws.Range("F42").Value = 41073477
$ws.Range("G42").Value = -1.74013

$ws.Range("B43").Value = "NEAR"
$ws.Range("C43").Value = "NEAR Protocol"
$ws.Range("D43").Value = 1.12
$ws.Range("E43").Value = 1050057796
$ws.Range("F43").Value = 28951221
$ws.Range("G43").Value = -0.92431

$ws.Range("B44").Value = "ARB"
$ws.Range("C44").Value = "Arbitrum"
$ws.Range("D44").Value = 0.814459
$ws.Range("E44").Value = 1038670651
$ws.Range("F44").Value = 80208761
$ws.Range("G44").Value = -1.29652

$ws.Range("B45").Value = "KAS"
$ws.Range("C45").Value = "Kaspa"
$ws.Range("D45").Value = 0.04594901
$ws.Range("E45").Value = 950618857
$ws.Range("F45").Value = 12786607
$ws.Range("G45").Value = -7.86616

$ws.Range("B46").Value = "RETH"
$ws.Range("C46").Value = "Rocket Pool ETH"
$ws.Range("D46").Value = 1770.79
$ws.Range("E46").Value = 932330934
$ws.Range("F46").Value = 2911711
$ws.Range("G46").Value = -0.24243

$ws.Range("D47").Value = 61.02
$ws.Range("E47").Value = 885940889
$ws.Range("F47").Value = 98662769
$ws.Range("G47").Value = 3.12624

$ws.Range("D48").Value = 0.086552
$ws.Range("E48").Value = 794266152
$ws.Range("F48").Value = 18397855
$ws.Range("G48").Value = 0.20492

$ws.Range("D49").Value = 5.28
$ws.Range("E49").Value = 760123446
$ws.Range("F49").Value = 6746252
$ws.Range("G49").Value = -0.98405

$ws.Range("B50").Value = "ALGO"
$ws.Range("C50").Value = "Algorand"
$ws.Range("D50").Value = 0.096162
$ws.Range("E50").Value = 752848302
$ws.Range("F50").Value = 15907026
$ws.Range("G50").Value = -0.29322

$ws.Range("B51").Value = "XDC"
$ws.Range("C51").Value = "XDC Network"
$ws.Range("D51").Value = 0.053273
$ws.Range("E51").Value = 738987447
$ws.Range("F51").Value = 5486661
$ws.Range("G51").Value = -2.0207

